# Adds new interview-experience rows to the log sheet.
#
# Before:
#   row13 = Starthealth (ht 240)
#   row14 = Accelya (ht 199.5, Result empty)
#
# After:
#   row13 = Here Technology            (new - chronologically earliest, inserted above)
#   row14 = Starthealth                (unchanged, just shifted down)
#   row15 = Accelya                    (unchanged, Result filled in = "cleared")
#   row16 = Accelya / Manager round    (new)
#   row17 = EY                         (new)
#   row18 = Accelya / Hr round         (new)
#   row19 = Clover infotech , andheri  (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the missing "Result" cell for the existing Accelya row (row 14,
#    becomes row 15 once the new row is inserted at the top later).
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = "cleared"

# ---------------------------------------------------------------------------
# 2) Append the new Accelya interview rounds (Manager round, EY walk-in,
#    Hr round) and the Clover infotech row directly below - rows 15-18 are
#    still unused at this point so no shifting is needed yet.
# ---------------------------------------------------------------------------

# Row 15 -> Accelya - Manager round - cleared
$ws.Range("A15").Value = 45891
$ws.Range("B15").Value = "Accelya"
$ws.Range("C15").Value = "Manager round"
$ws.Range("D15").Value = "cleared"

# Row 16 -> EY - collections, elastic search, oops,  - cleared
$ws.Range("A16").Value = 45891
$ws.Range("B16").Value = "EY"
$ws.Range("C16").Value = "collections, elastic search, oops, "
$ws.Range("C16").VerticalAlignment = -4160
$ws.Range("D16").Value = "cleared"
$ws.Rows.Item(16).RowHeight = 29.25

# Row 17 -> Accelya - Hr round - cleared. Got offer
$ws.Range("A17").Value = 45894
$ws.Range("B17").Value = "Accelya"
$ws.Range("C17").Value = "Hr round"
$ws.Range("D17").Value = "cleared. Got offer"

# Row 18 -> Clover infotech , andheri
$ws.Range("A18").Value = 45914
$ws.Range("B18").Value = "Clover infotech , andheri"
$ws.Range("C18").Value = "java - callable vs runnable, java 8 features, checked vs unchecked exceptions, how to implement thread, default and static method, internal working of hashmap, fail fast  and fail safe iterator, code - print frequecy of each vowel in string" + [char]10 + "spring boot - global exception, spring security"
$ws.Range("D18").Value = "cleared"
$ws.Rows.Item(18).RowHeight = 90

# ---------------------------------------------------------------------------
# 3) Insert the older "Here Technology" interview at its chronological spot,
#    row 13 (it happened on 2025-07-31, before the Starthealth interview),
#    pushing every row below it down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A14:D14").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A13").Value = 45869
$ws.Range("B13").Value = "Here rechnology"
$ws.Range("C13").Value = "git, write unit tests , linux commands, project , pagination"
$ws.Range("D13").Value = "failed"

# ---------------------------------------------------------------------------
# 4) Leave the view/selection the way it was left after the edit.
# ---------------------------------------------------------------------------
$null = $ws.Range("D14").Select()

Write-Output "done"
